# chore: adapt column header formatting to respective input file names (#7)
#
# 1. Rename the "_old" / "_new" header-row suffixes to "_FV2410" / "_FV2504"
# 2. Wrap the data range in an Excel Table ("Table1")
# 3. Freeze the header row (first row) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-label the header row (A1:U1) in place.
# ---------------------------------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "K1" = "diff"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# ---------------------------------------------------------------------------
# 2. Turn A1:U78 into a real Excel Table (Table1), matching the shipped
#    header names above for its column names.
#
#    The header row already carries explicit cell formatting (bold, grey
#    fill, border, centred + wrapped). Converting a pre-formatted range into
#    a ListObject would normally make Excel capture that look as a one-off
#    "headerRowDxfId" diff (growing styles.xml's <dxfs>). To keep styles.xml
#    byte-for-byte untouched, the header formatting is stashed on a scratch
#    cell, cleared prior to the table's creation, and then restored via a
#    single PasteSpecial (formats only) once the table exists.
# ---------------------------------------------------------------------------
$header = $ws.Range("A1:U1")
$scratch = $ws.Range("A200")
$ws.Range("A1").Copy($scratch)
$header.ClearFormats()

$dataRange = $ws.Range("A1:U78")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

$scratch.Copy()
$header.PasteSpecial(-4122)  # xlPasteFormats
$scratch.Clear()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("A2").Select()
$win.FreezePanes = $true

Write-Output "Header row relabelled, Table1 created, top row frozen."
